# Fruta / hortaliza, semanal
# Insert two new weekly price rows (row 182 and 183) just above the existing
# data block that starts at row 182, shifting the previous rows 182-211 down
# to 184-213. The new rows carry a fresh "Especial"/"Primera" quality pair
# for the Packham's Triumph pear sold at Terminal Hortofruticola Agro
# Chillan, dated 2022-08-03 (serial 44776).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 182:211 down to 184:213, creating two blank rows.
$ws.Rows("182:183").Insert()

# ---- New row 182 : Especial ----------------------------------------------
$ws.Cells.Item(182, 1).Value = 7
$ws.Cells.Item(182, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(182, 3).Value = "Ñuble"
$ws.Cells.Item(182, 4).Value = 44776
$ws.Cells.Item(182, 5).Value = 16
$ws.Cells.Item(182, 6).Value = "Fruta"
$ws.Cells.Item(182, 7).Value = 100104
$ws.Cells.Item(182, 8).Value = "Frutos de pepita"
$ws.Cells.Item(182, 9).Value = 100104005
$ws.Cells.Item(182, 10).Value = "Pera"
$ws.Cells.Item(182, 11).Value = "Packham's Triumph"
$ws.Cells.Item(182, 12).Value = "Especial"
$ws.Cells.Item(182, 13).Value = 60
$ws.Cells.Item(182, 14).Value = 10000
$ws.Cells.Item(182, 15).Value = 10000
$ws.Cells.Item(182, 16).Value = 10000
$ws.Cells.Item(182, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(182, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(182, 19).Value = 625
$ws.Cells.Item(182, 20).Value = 16

# ---- New row 183 : Primera ------------------------------------------------
$ws.Cells.Item(183, 1).Value = 7
$ws.Cells.Item(183, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(183, 3).Value = "Ñuble"
$ws.Cells.Item(183, 4).Value = 44776
$ws.Cells.Item(183, 5).Value = 16
$ws.Cells.Item(183, 6).Value = "Fruta"
$ws.Cells.Item(183, 7).Value = 100104
$ws.Cells.Item(183, 8).Value = "Frutos de pepita"
$ws.Cells.Item(183, 9).Value = 100104005
$ws.Cells.Item(183, 10).Value = "Pera"
$ws.Cells.Item(183, 11).Value = "Packham's Triumph"
$ws.Cells.Item(183, 12).Value = "Primera"
$ws.Cells.Item(183, 13).Value = 120
$ws.Cells.Item(183, 14).Value = 8500
$ws.Cells.Item(183, 15).Value = 9000
$ws.Cells.Item(183, 16).Value = 8750
$ws.Cells.Item(183, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(183, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(183, 19).Value = 547
$ws.Cells.Item(183, 20).Value = 16
